$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 459.5
$ws.Range("I8").Value = 459.5
$ws.Range("K8").Value = 1378.5
$ws.Range("M8").Value = -1239.5
$ws.Range("H39").Value = 495.27274
$ws.Range("I39").Value = 393.5
$ws.Range("J39").Value = 766.6667
$ws.Range("K39").Value = 1180.5
$ws.Range("L39").Value = 2300.0001
$ws.Range("M39").Value = -884.5
$ws.Range("N39").Value = -2892.0001
$ws.Range("H98").Value = 2886417.2
$ws.Range("I98").Value = 3136962.8
$ws.Range("K98").Value = 3136962.8
$ws.Range("M98").Value = -3135464.8
$ws.Range("H107").Value = 15627458
$ws.Range("I107").Value = 9617895
$ws.Range("J107").Value = 41668896
$ws.Range("K107").Value = 9617895
$ws.Range("L107").Value = 41668896
$ws.Range("M107").Value = -9615975
$ws.Range("N107").Value = -41672736
$ws.Range("H113").Value = 6690.75
$ws.Range("I113").Value = 10135.467
$ws.Range("J113").Value = 3651.2942
$ws.Range("K113").Value = 10135.467
$ws.Range("L113").Value = 3651.2942
$ws.Range("M113").Value = -6881.467000000001
$ws.Range("N113").Value = -10159.2942
$ws.Range("H116").Value = 3673.125
$ws.Range("I116").Value = 3653.5557
$ws.Range("K116").Value = 3653.5557
$ws.Range("M116").Value = -211.5556999999999
$ws.Range("H122").Value = 2886417.2
$ws.Range("I122").Value = 3136962.8
$ws.Range("K122").Value = 9410888.399999999
$ws.Range("M122").Value = -9408438.399999999
$ws.Range("H135").Value = 1227.9474
$ws.Range("I135").Value = 830.5294
$ws.Range("K135").Value = 7474.7646
$ws.Range("M135").Value = -4939.7646
$ws.Range("H138").Value = 3984.8936
$ws.Range("I138").Value = 3234.8572
$ws.Range("J138").Value = 4045.2415
$ws.Range("K138").Value = 9704.571599999999
$ws.Range("L138").Value = 12135.7245
$ws.Range("M138").Value = -4564.571599999999
$ws.Range("N138").Value = -22415.7245
$ws.Range("H141").Value = 1491.5294
$ws.Range("I141").Value = 1491.5294
$ws.Range("K141").Value = 4474.5882
$ws.Range("M141").Value = 705.4117999999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1806.25
$ws.Range("I45").Value = 1703.8
$ws.Range("K45").Value = 1703.8
$ws.Range("M45").Value = -1326.8
$ws.Range("H74").Value = 93018.82000000001
$ws.Range("I74").Value = 118558.766
$ws.Range("J74").Value = 6183
$ws.Range("K74").Value = 118558.766
$ws.Range("L74").Value = 6183
$ws.Range("M74").Value = -117684.766
$ws.Range("N74").Value = -7931
$ws.Range("H77").Value = 93018.82000000001
$ws.Range("I77").Value = 118558.766
$ws.Range("J77").Value = 6183
$ws.Range("K77").Value = 592793.8300000001
$ws.Range("L77").Value = 30915
$ws.Range("M77").Value = -588425.8300000001
$ws.Range("N77").Value = -39651
$ws.Range("H110").Value = 6806.1333
$ws.Range("I110").Value = 2917.5454
$ws.Range("K110").Value = 2917.5454
$ws.Range("M110").Value = -872.5454
$ws.Range("H122").Value = 2228.36
$ws.Range("I122").Value = 2189.5417
$ws.Range("J122").Value = 3160
$ws.Range("K122").Value = 6568.625100000001
$ws.Range("L122").Value = 9480
$ws.Range("M122").Value = -4118.625100000001
$ws.Range("N122").Value = -14380
$ws.Range("H123").Value = 67000
$ws.Range("I123").Value = 67000
$ws.Range("K123").Value = 67000
$ws.Range("M123").Value = -62100
$ws.Range("H132").Value = 3500.578
$ws.Range("I132").Value = 3062.3901
$ws.Range("J132").Value = 7992
$ws.Range("K132").Value = 9187.1703
$ws.Range("L132").Value = 23976
$ws.Range("M132").Value = -6657.1703
$ws.Range("N132").Value = -29036
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2109.0789
$ws.Range("I20").Value = 1445.9286
$ws.Range("K20").Value = 1445.9286
$ws.Range("M20").Value = -1198.9286
$ws.Range("H105").Value = 2935.2
$ws.Range("I105").Value = 2935.2
$ws.Range("K105").Value = 2935.2
$ws.Range("M105").Value = -1188.2
$ws.Range("H134").Value = 2710.457
$ws.Range("I134").Value = 2159.742
$ws.Range("J134").Value = 6978.5
$ws.Range("K134").Value = 6479.226000000001
$ws.Range("L134").Value = 20935.5
$ws.Range("M134").Value = -3944.226000000001
$ws.Range("N134").Value = -26005.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 2828.2942
$ws.Range("I12").Value = 2720.6428
$ws.Range("J12").Value = 3330.6667
$ws.Range("K12").Value = 2720.6428
$ws.Range("L12").Value = 3330.6667
$ws.Range("M12").Value = -2550.6428
$ws.Range("N12").Value = -3670.6667
$ws.Range("H58").Value = 1482.6216
$ws.Range("I58").Value = 1386.1818
$ws.Range("J58").Value = 2278.25
$ws.Range("K58").Value = 1386.1818
$ws.Range("L58").Value = 2278.25
$ws.Range("M58").Value = -1183.1818
$ws.Range("N58").Value = -2684.25
$ws.Range("H122").Value = 2598.7646
$ws.Range("I122").Value = 2180.75
$ws.Range("K122").Value = 6542.25
$ws.Range("M122").Value = -4092.25
$ws.Range("H136").Value = 1482.6216
$ws.Range("I136").Value = 1386.1818
$ws.Range("J136").Value = 2278.25
$ws.Range("K136").Value = 4158.5454
$ws.Range("L136").Value = 6834.75
$ws.Range("M136").Value = -1608.5454
$ws.Range("N136").Value = -11934.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 5160.5
$ws.Range("J132").Value = 2252.111
$ws.Range("L132").Value = 20268.999
$ws.Range("N132").Value = -25328.999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8643.817999999999
$ws.Range("I70").Value = 9038
$ws.Range("J70").Value = 7592.6665
$ws.Range("K70").Value = 9038
$ws.Range("L70").Value = 7592.6665
$ws.Range("M70").Value = -8768
$ws.Range("N70").Value = -8132.6665
$ws.Range("H73").Value = 8643.817999999999
$ws.Range("I73").Value = 9038
$ws.Range("J73").Value = 7592.6665
$ws.Range("K73").Value = 9038
$ws.Range("L73").Value = 7592.6665
$ws.Range("M73").Value = -8102
$ws.Range("N73").Value = -9464.666499999999
$ws.Range("H80").Value = 4176.8
$ws.Range("I80").Value = 3176.7778
$ws.Range("K80").Value = 3176.7778
$ws.Range("M80").Value = -2178.7778
$ws.Range("H83").Value = 4176.8
$ws.Range("I83").Value = 3176.7778
$ws.Range("K83").Value = 15883.889
$ws.Range("M83").Value = -10891.889
$ws.Range("H92").Value = 7333.3335
$ws.Range("J92").Value = 7333.3335
$ws.Range("L92").Value = 7333.3335
$ws.Range("N92").Value = -11077.3335
$ws.Range("H113").Value = 9598
$ws.Range("I113").Value = 2998
$ws.Range("J113").Value = 19498
$ws.Range("K113").Value = 2998
$ws.Range("L113").Value = 19498
$ws.Range("M113").Value = -828
$ws.Range("N113").Value = -23838
$ws.Range("H122").Value = 1193156.9
$ws.Range("I122").Value = 1925882.6
$ws.Range("K122").Value = 5777647.800000001
$ws.Range("M122").Value = -5775197.800000001
$ws.Range("H132").Value = 60128.79
$ws.Range("I132").Value = 60128.79
$ws.Range("K132").Value = 180386.37
$ws.Range("M132").Value = -177856.37
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 47749.125
$ws.Range("J38").Value = 47833.168
$ws.Range("L38").Value = 47833.168
$ws.Range("N38").Value = -48653.168
$ws.Range("H59").Value = 27500
$ws.Range("J59").Value = 27500
$ws.Range("L59").Value = 27500
$ws.Range("N59").Value = -28808
$ws.Range("H122").Value = 5854.875
$ws.Range("I122").Value = 5844.4614
$ws.Range("J122").Value = 5900
$ws.Range("K122").Value = 17533.3842
$ws.Range("L122").Value = 17700
$ws.Range("M122").Value = -15083.3842
$ws.Range("N122").Value = -22600
$ws.Range("H132").Value = 3760.3867
$ws.Range("I132").Value = 2712.8293
$ws.Range("J132").Value = 5023.6177
$ws.Range("K132").Value = 8138.4879
$ws.Range("L132").Value = 15070.8531
$ws.Range("M132").Value = -5608.4879
$ws.Range("N132").Value = -20130.8531
$ws.Range("H136").Value = 3488.0881
$ws.Range("I136").Value = 1996
$ws.Range("K136").Value = 5988
$ws.Range("M136").Value = -3438
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12246.4
$ws.Range("I62").Value = 9997
$ws.Range("K62").Value = 9997
$ws.Range("M62").Value = -9373
$ws.Range("H65").Value = 12246.4
$ws.Range("I65").Value = 9997
$ws.Range("K65").Value = 49985
$ws.Range("M65").Value = -46865
$ws.Range("H96").Value = 88535.086
$ws.Range("I96").Value = 171605.33
$ws.Range("J96").Value = 5464.8335
$ws.Range("K96").Value = 171605.33
$ws.Range("L96").Value = 5464.8335
$ws.Range("M96").Value = -170232.33
$ws.Range("N96").Value = -8210.833500000001
$ws.Range("H122").Value = 2469.5789
$ws.Range("I122").Value = 2195.2
$ws.Range("J122").Value = 3498.5
$ws.Range("K122").Value = 6585.599999999999
$ws.Range("L122").Value = 10495.5
$ws.Range("M122").Value = -4135.599999999999
$ws.Range("N122").Value = -15395.5
$ws.Range("H136").Value = 440517.53
$ws.Range("I136").Value = 506025.4
$ws.Range("J136").Value = 3798.3333
$ws.Range("K136").Value = 1518076.2
$ws.Range("L136").Value = 11394.9999
$ws.Range("M136").Value = -1515526.2
$ws.Range("N136").Value = -16494.9999
